$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 360, shifting existing rows 360:382 down to 361:383
$ws.Rows(360).Insert()

# Populate the newly inserted row with the new weekly record
$ws.Range("A360").Value = 9
$ws.Range("B360").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C360").Value = "Metropolitana"
$ws.Range("D360").Value = 44706
$ws.Range("E360").Value = 13
$ws.Range("F360").Value = 100112052
$ws.Range("G360").Value = "Albahaca"
$ws.Range("H360").Value = "Sin especificar"
$ws.Range("I360").Value = "Primera"
$ws.Range("J360").Value = 230
$ws.Range("K360").Value = 4000
$ws.Range("L360").Value = 5000
$ws.Range("M360").Value = 4435
$ws.Range("N360").Value = '$/paquete'
$ws.Range("O360").Value = "Provincia de Chacabuco"
$ws.Range("P360").Value = 4435
$ws.Range("Q360").Value = 1
$ws.Range("R360").Value = "Hortaliza"
